# Generate Report for Handoff
#
# A fresh localization handoff run produced a new source-file GUID and a
# new content hash for the generated xliff files; refresh the status
# report with the new file names and timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "c6c8f945-f02c-442e-b7a2-bdf442de0434"
$newGuid = "2b4f149c-ea37-4f1f-ac12-f9cd1fd61d1f"

$oldHash = "bb6b41f172739fbb25c77a565c05115385e05c7f"
$newHash = "715ac7f8c2e9232e81c2b3ee19f6ad47a4261152"

$newHoDate = "2016-08-27 20:56:07"
$newZhDate = "2016-08-27 20:55:59"

# The external hyperlink target itself (GitHub blob URL) is unchanged by
# this handoff run -- only the on-sheet display text is refreshed.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68ae05cff7e7bc2954bbe0fe37073c518f4d9791/e2e/$oldGuid.md"

function Update-Hyperlink($ws, $cellRef, $displayText) {
    $rng = $ws.Range($cellRef)
    $rng.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($rng, $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $displayText) | Out-Null
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
Update-Hyperlink $wsOverview "B2" "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newHoDate

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-Hyperlink $wsZh "A2" "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $newZhDate

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
Update-Hyperlink $wsDe "A2" "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = $newHoDate

# --- Column A width tweak on all three sheets ---
# Target stored width is 39.3234034946987 characters; the host's column
# width setter snaps to the nearest 1/6-character increment, so feed it
# the input that lands on the closest representable width (39.333...).
$wsOverview.Columns.Item(1).ColumnWidth = 38.5
$wsZh.Columns.Item(1).ColumnWidth = 38.5
$wsDe.Columns.Item(1).ColumnWidth = 38.5
